$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$row = 57

# Columns A (date-like "2025-02-07") and D ("05") look numeric/date to Excel's
# smart-entry parser, so force them to Text before writing, then restore the
# default "General" format so the saved cell carries no explicit style.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-02-07"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = "22:32:30"
$ws.Cells.Item($row, 3).Value = "Friday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "05"
$ws.Cells.Item($row, 4).Style = "Normal"

$ws.Cells.Item($row, 5).Value = 126230
$ws.Cells.Item($row, 6).Value = 141779
$ws.Cells.Item($row, 7).Value = 168193
$ws.Cells.Item($row, 8).Value = 158309
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 143246
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 191539
$ws.Cells.Item($row, 14).Value = 115195
$ws.Cells.Item($row, 15).Value = 44781
$ws.Cells.Item($row, 16).Value = 28317
$ws.Cells.Item($row, 17).Value = 63887
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 41739
$ws.Cells.Item($row, 20).Value = -1
